# Auto-generated update of FFXIV Leve profit-tracking values
# (scheduled-runner refresh of Universalis market-price snapshots).
# For each sheet, a table of cell -> new numeric value is applied
# via Range.Value, matching the refreshed currentAveragePrice /
# LevePrice / LeveProfit columns (H, I, J, K, L, M, N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @("H32", 6122.2),
    @("J32", 5499.25),
    @("L32", 5499.25),
    @("N32", -6151.25),
    @("H48", 5128.5713),
    @("J48", 5128.5713),
    @("L48", 15385.7139),
    @("N48", -15969.7139),
    @("H56", 5128.5713),
    @("J56", 5128.5713),
    @("L56", 15385.7139),
    @("N56", -16453.7139),
    @("H113", 61287),
    @("I113", 75940.86),
    @("K113", 75940.86),
    @("M113", -72686.86),
    @("H132", 3690.0645),
    @("I132", 3744.6296),
    @("K132", 11233.8888),
    @("M132", -8703.888800000001),
    @("H138", 2102.9443),
    @("I138", 1593.1111),
    @("J138", 2612.7778),
    @("K138", 4779.3333),
    @("L138", 7838.3334),
    @("M138", 360.6666999999998),
    @("N138", -18118.3334),
    @("H139", 98255.60000000001),
    @("J139", 98255.60000000001),
    @("L139", 98255.60000000001),
    @("N139", -108535.6)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @("H102", 590415.0600000001),
    @("I102", 716544),
    @("K102", 716544),
    @("M102", -714922),
    @("H122", 14494098),
    @("I122", 1473.1111),
    @("J122", 66667548),
    @("K122", 4419.3333),
    @("L122", 200002644),
    @("M122", -1969.3333),
    @("N122", -200007544),
    @("H132", 1697.6938),
    @("I132", 1660.3096),
    @("K132", 4980.9288),
    @("M132", -2450.9288),
    @("H137", 89743.89999999999),
    @("J137", 89743.89999999999),
    @("L137", 89743.89999999999),
    @("N137", -99943.89999999999),
    @("H140", 69947.5),
    @("J140", 69947.5),
    @("L140", 69947.5),
    @("N140", -80307.5)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @("H22", 915.5),
    @("J22", 900.5),
    @("L22", 900.5),
    @("N22", -1246.5),
    @("H60", 35864.168),
    @("I60", 7999),
    @("J60", 41437.2),
    @("K60", 7999),
    @("L60", 41437.2),
    @("M60", -7400),
    @("N60", -42635.2),
    @("H81", 24470.25),
    @("J81", 24470.25),
    @("L81", 24470.25),
    @("N81", -26592.25),
    @("H84", 24470.25),
    @("J84", 24470.25),
    @("L84", 73410.75),
    @("N84", -84018.75),
    @("H86", 58825780),
    @("I86", 83335310),
    @("K86", 83335310),
    @("M86", -83334187),
    @("H89", 58825780),
    @("I89", 83335310),
    @("K89", 416676550),
    @("M89", -416670934),
    @("H99", 1649.7931),
    @("I99", 1612.5),
    @("J99", 1828.8),
    @("K99", 1612.5),
    @("L99", 1828.8),
    @("M99", -114.5),
    @("N99", -4824.8),
    @("H133", 79740),
    @("J133", 79740),
    @("L133", 79740),
    @("N133", -89860),
    @("H138", 77687.71000000001),
    @("J138", 77687.71000000001),
    @("L138", 77687.71000000001),
    @("N138", -87967.71000000001),
    @("H140", 72984.10000000001),
    @("J140", 72984.10000000001),
    @("L140", 72984.10000000001),
    @("N140", -83344.10000000001)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @("H52", 115470.75),
    @("J52", 115470.75),
    @("L52", 115470.75),
    @("N52", -116058.75),
    @("H58", 46669530),
    @("I58", 33335666),
    @("J58", 55558772),
    @("K58", 33335666),
    @("L58", 55558772),
    @("M58", -33335463),
    @("N58", -55559178),
    @("H107", 71431700),
    @("I107", 125001790),
    @("J107", 4892.6665),
    @("K107", 125001790),
    @("L107", 4892.6665),
    @("M107", -124999870),
    @("N107", -8732.666499999999),
    @("H122", 762.6667),
    @("I122", 760.25),
    @("J122", 767.5),
    @("K122", 2280.75),
    @("L122", 2302.5),
    @("M122", 169.25),
    @("N122", -7202.5),
    @("H132", 1882.4642),
    @("I132", 1894.5),
    @("J132", 1810.25),
    @("K132", 5683.5),
    @("L132", 5430.75),
    @("M132", -3153.5),
    @("N132", -10490.75),
    @("H134", 2123.5),
    @("I134", 2125.7896),
    @("K134", 6377.3688),
    @("M134", -3842.3688),
    @("H136", 46669530),
    @("I136", 33335666),
    @("J136", 55558772),
    @("K136", 100006998),
    @("L136", 166676316),
    @("M136", -100004448),
    @("N136", -166681416),
    @("H139", 82376.8),
    @("J139", 90293.75),
    @("L139", 90293.75),
    @("N139", -100573.75),
    @("H140", 69923.625),
    @("J140", 69923.625),
    @("L140", 69923.625),
    @("N140", -80283.625)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @("H102", 412908.44),
    @("I102", 469503.62),
    @("J102", 7309.6665),
    @("K102", 469503.62),
    @("L102", 7309.6665),
    @("M102", -467881.62),
    @("N102", -10553.6665),
    @("H122", 100502.37),
    @("J122", 134998),
    @("L122", 404994),
    @("N122", -409894),
    @("H132", 8954.454),
    @("I132", 8449.9),
    @("K132", 25349.7),
    @("M132", -22819.7),
    @("H135", 86652.5),
    @("J135", 86652.5),
    @("L135", 86652.5),
    @("N135", -96792.5)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @("H22", 3773.111),
    @("I22", 3769),
    @("J22", 3778.25),
    @("K22", 3769),
    @("L22", 3778.25),
    @("M22", -3474),
    @("N22", -4368.25),
    @("H27", 3773.111),
    @("I27", 3769),
    @("J27", 3778.25),
    @("K27", 3769),
    @("L27", 3778.25),
    @("M27", -3662),
    @("N27", -3992.25),
    @("H55", 2515.875),
    @("I55", 1512.8572),
    @("K55", 1512.8572),
    @("M55", -1339.8572),
    @("H93", 856.7646999999999),
    @("I93", 746.04346),
    @("J93", 1088.2727),
    @("K93", 746.04346),
    @("L93", 1088.2727),
    @("M93", 501.95654),
    @("N93", -3584.2727),
    @("H100", 4646.1665),
    @("I100", 4828.1),
    @("J100", 3736.5),
    @("K100", 4828.1),
    @("L100", 3736.5),
    @("M100", -4287.1),
    @("N100", -4818.5),
    @("H122", 4591.75),
    @("I122", 3320.1667),
    @("J122", 6499.125),
    @("K122", 9960.500100000001),
    @("L122", 19497.375),
    @("M122", -7510.500100000001),
    @("N122", -24397.375),
    @("H132", 3777.2778),
    @("I132", 3366.1333),
    @("K132", 10098.3999),
    @("M132", -7568.3999)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @("H113", 1103.1666),
    @("I113", 1059.0476),
    @("J113", 1206.1111),
    @("K113", 3177.142800000001),
    @("L113", 3618.3333),
    @("M113", -1007.142800000001),
    @("N113", -7958.3333),
    @("H122", 2738.1035),
    @("I122", 2260.8948),
    @("K122", 6782.6844),
    @("M122", -4332.6844),
    @("H126", 83335460),
    @("I126", 166668580),
    @("K126", 500005740),
    @("M126", -500003270),
    @("H132", 983.55554),
    @("I132", 983.55554),
    @("K132", 2950.66662),
    @("M132", -420.66662)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

Write-Output "Applied 232 cell updates across 7 sheets"
